$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$g = $s.Shapes.Item(2)

# pl6 (id=7)
$sh = $g.GroupItems.Item(4)
$sh.Top = 472.68695068359375

# pl7 (id=8)
$sh = $g.GroupItems.Item(5)
$sh.Top = 339.1060791015625

# pl8 (id=9)
$sh = $g.GroupItems.Item(6)
$sh.Top = 205.52520751953125

# pl9 (id=10)
$sh = $g.GroupItems.Item(7)
$sh.Top = 539.4773559570312

# pl10 (id=11)
$sh = $g.GroupItems.Item(8)
$sh.Top = 405.896484375

# pl11 (id=12)
$sh = $g.GroupItems.Item(9)
$sh.Top = 272.31561279296875

# pl12 (id=13)
$sh = $g.GroupItems.Item(10)
$sh.Top = 138.7348175048828

# pg15 (id=16)
$sh = $g.GroupItems.Item(13)
$sh.Left = 155.86441040039062
$sh.Top = 259.9656982421875

# pg16 (id=17)
$sh = $g.GroupItems.Item(14)
$sh.Left = 195.5236358642578
$sh.Top = 259.8893127441406

# pg17 (id=18)
$sh = $g.GroupItems.Item(15)
$sh.Left = 180.3412628173828
$sh.Top = 306.72900390625

# pg18 (id=19)
$sh = $g.GroupItems.Item(16)
$sh.Left = 167.83921813964844
$sh.Top = 261.2213439941406

# pg19 (id=20)
$sh = $g.GroupItems.Item(17)
$sh.Left = 134.0924530029297
$sh.Top = 295.47308349609375

# pg20 (id=21)
$sh = $g.GroupItems.Item(18)
$sh.Left = 174.70811462402344
$sh.Top = 291.5550537109375

# pl23 (id=24)
$sh = $g.GroupItems.Item(21)
$sh.Top = 472.68695068359375

# pl24 (id=25)
$sh = $g.GroupItems.Item(22)
$sh.Top = 339.1060791015625

# pl25 (id=26)
$sh = $g.GroupItems.Item(23)
$sh.Top = 205.52520751953125

# pl26 (id=27)
$sh = $g.GroupItems.Item(24)
$sh.Top = 539.4773559570312

# pl27 (id=28)
$sh = $g.GroupItems.Item(25)
$sh.Top = 405.896484375

# pl28 (id=29)
$sh = $g.GroupItems.Item(26)
$sh.Top = 272.31561279296875

# pl29 (id=30)
$sh = $g.GroupItems.Item(27)
$sh.Top = 138.7348175048828

# pg32 (id=33)
$sh = $g.GroupItems.Item(30)
$sh.Left = 533.0557861328125
$sh.Top = 473.9181213378906

# pg33 (id=34)
$sh = $g.GroupItems.Item(31)
$sh.Left = 468.7443542480469
$sh.Top = 447.7611999511719

# pl36 (id=37)
$sh = $g.GroupItems.Item(34)
$sh.Top = 472.68695068359375

# pl37 (id=38)
$sh = $g.GroupItems.Item(35)
$sh.Top = 339.1060791015625

# pl38 (id=39)
$sh = $g.GroupItems.Item(36)
$sh.Top = 205.52520751953125

# pl39 (id=40)
$sh = $g.GroupItems.Item(37)
$sh.Top = 539.4773559570312

# pl40 (id=41)
$sh = $g.GroupItems.Item(38)
$sh.Top = 405.896484375

# pl41 (id=42)
$sh = $g.GroupItems.Item(39)
$sh.Top = 272.31561279296875

# pl42 (id=43)
$sh = $g.GroupItems.Item(40)
$sh.Top = 138.7348175048828

# pg45 (id=46)
$sh = $g.GroupItems.Item(43)
$sh.Left = 654.0865478515625
$sh.Top = 234.4820556640625

# pg46 (id=47)
$sh = $g.GroupItems.Item(44)
$sh.Left = 651.267578125
$sh.Top = 255.7223663330078

# pg47 (id=48)
$sh = $g.GroupItems.Item(45)
$sh.Left = 762.97119140625

# pg48 (id=49)
$sh = $g.GroupItems.Item(46)
$sh.Left = 710.8084716796875
$sh.Top = 508.90118408203125

# pg49 (id=50)
$sh = $g.GroupItems.Item(47)
$sh.Left = 705.2887573242188
$sh.Top = 406.2446594238281

# pg50 (id=51)
$sh = $g.GroupItems.Item(48)
$sh.Left = 586.2157592773438
$sh.Top = 269.4608154296875

# pg51 (id=52)
$sh = $g.GroupItems.Item(49)
$sh.Left = 724.8751220703125
$sh.Top = 483.8338623046875

# pg52 (id=53)
$sh = $g.GroupItems.Item(50)
$sh.Left = 758.5652465820312
$sh.Top = 455.20782470703125

# pg53 (id=54)
$sh = $g.GroupItems.Item(51)
$sh.Left = 621.5642700195312
$sh.Top = 366.37811279296875

# pg54 (id=55)
$sh = $g.GroupItems.Item(52)
$sh.Left = 592.1602172851562
$sh.Top = 179.52252197265625

# pg55 (id=56)
$sh = $g.GroupItems.Item(53)
$sh.Left = 653.9293212890625
$sh.Top = 226.79481506347656

# pg56 (id=57)
$sh = $g.GroupItems.Item(54)
$sh.Left = 606.5917358398438
$sh.Top = 194.17835998535156

# pg57 (id=58)
$sh = $g.GroupItems.Item(55)
$sh.Left = 608.839111328125
$sh.Top = 371.3288269042969

# pg58 (id=59)
$sh = $g.GroupItems.Item(56)
$sh.Left = 628.9163818359375
$sh.Top = 202.99732971191406

# pg59 (id=60)
$sh = $g.GroupItems.Item(57)
$sh.Left = 591.684814453125
$sh.Top = 176.23985290527344

# pg60 (id=61)
$sh = $g.GroupItems.Item(58)
$sh.Left = 644.3438720703125
$sh.Top = 176.2911834716797

# pg61 (id=62)
$sh = $g.GroupItems.Item(59)
$sh.Left = 649.1234130859375
$sh.Top = 176.238037109375

# pg62 (id=63)
$sh = $g.GroupItems.Item(60)
$sh.Left = 656.5753784179688
$sh.Top = 176.29678344726562

# pg63 (id=64)
$sh = $g.GroupItems.Item(61)
$sh.Left = 637.8712158203125
$sh.Top = 176.3002471923828

# pg64 (id=65)
$sh = $g.GroupItems.Item(62)
$sh.Left = 610.2740478515625
$sh.Top = 380.8790588378906

# pg65 (id=66)
$sh = $g.GroupItems.Item(63)
$sh.Left = 660.1895751953125
$sh.Top = 353.95294189453125

# pg66 (id=67)
$sh = $g.GroupItems.Item(64)
$sh.Left = 590.27099609375
$sh.Top = 203.09835815429688

# pg67 (id=68)
$sh = $g.GroupItems.Item(65)
$sh.Left = 601.4033203125
$sh.Top = 263.8846740722656

# pg68 (id=69)
$sh = $g.GroupItems.Item(66)
$sh.Left = 638.685302734375
$sh.Top = 252.38433837890625

# pg69 (id=70)
$sh = $g.GroupItems.Item(67)
$sh.Left = 594.7350463867188
$sh.Top = 388.55401611328125

# pg70 (id=71)
$sh = $g.GroupItems.Item(68)
$sh.Left = 592.8582763671875
$sh.Top = 358.4952087402344

# pg71 (id=72)
$sh = $g.GroupItems.Item(69)
$sh.Left = 763.7293090820312
$sh.Top = 402.80120849609375

# pg72 (id=73)
$sh = $g.GroupItems.Item(70)
$sh.Left = 652.8309936523438
$sh.Top = 160.61119079589844

# pg73 (id=74)
$sh = $g.GroupItems.Item(71)
$sh.Left = 582.71484375
$sh.Top = 246.2820587158203

# pg74 (id=75)
$sh = $g.GroupItems.Item(72)
$sh.Left = 595.852294921875
$sh.Top = 270.9327697753906

# pg75 (id=76)
$sh = $g.GroupItems.Item(73)
$sh.Left = 612.56689453125
$sh.Top = 354.7245178222656

# pg76 (id=77)
$sh = $g.GroupItems.Item(74)
$sh.Left = 608.5021362304688
$sh.Top = 189.7017364501953

# pg77 (id=78)
$sh = $g.GroupItems.Item(75)
$sh.Left = 661.5023193359375
$sh.Top = 179.30804443359375

# pg78 (id=79)
$sh = $g.GroupItems.Item(76)
$sh.Left = 627.7211303710938
$sh.Top = 280.534912109375

# pg79 (id=80)
$sh = $g.GroupItems.Item(77)
$sh.Left = 658.2767944335938
$sh.Top = 356.06536865234375

# pg80 (id=81)
$sh = $g.GroupItems.Item(78)
$sh.Left = 601.840087890625
$sh.Top = 170.1629180908203

# pg81 (id=82)
$sh = $g.GroupItems.Item(79)
$sh.Left = 604.1676025390625
$sh.Top = 194.675048828125

# pg82 (id=83)
$sh = $g.GroupItems.Item(80)
$sh.Left = 649.3206787109375
$sh.Top = 164.25811767578125

# pg83 (id=84)
$sh = $g.GroupItems.Item(81)
$sh.Left = 636.0950927734375
$sh.Top = 194.45425415039062

# pg84 (id=85)
$sh = $g.GroupItems.Item(82)
$sh.Left = 636.9769897460938
$sh.Top = 160.32386779785156

# pg85 (id=86)
$sh = $g.GroupItems.Item(83)
$sh.Left = 600.4698486328125
$sh.Top = 216.18040466308594

# pg86 (id=87)
$sh = $g.GroupItems.Item(84)
$sh.Left = 652.9053955078125

# pg87 (id=88)
$sh = $g.GroupItems.Item(85)
$sh.Left = 607.2913818359375
$sh.Top = 254.54229736328125

# pg88 (id=89)
$sh = $g.GroupItems.Item(86)
$sh.Left = 614.8023071289062
$sh.Top = 207.11590576171875

# pg89 (id=90)
$sh = $g.GroupItems.Item(87)
$sh.Left = 640.1900024414062
$sh.Top = 197.77496337890625

# pg90 (id=91)
$sh = $g.GroupItems.Item(88)
$sh.Left = 624.79150390625
$sh.Top = 214.8120574951172

# pg91 (id=92)
$sh = $g.GroupItems.Item(89)
$sh.Left = 595.010498046875
$sh.Top = 253.5963134765625

# pg92 (id=93)
$sh = $g.GroupItems.Item(90)
$sh.Left = 616.8720092773438
$sh.Top = 231.8785858154297

# pg93 (id=94)
$sh = $g.GroupItems.Item(91)
$sh.Left = 598.2661743164062
$sh.Top = 243.672607421875

# pg94 (id=95)
$sh = $g.GroupItems.Item(92)
$sh.Left = 619.3300170898438
$sh.Top = 236.98985290527344

# pg95 (id=96)
$sh = $g.GroupItems.Item(93)
$sh.Left = 647.09326171875
$sh.Top = 259.4040222167969

# pg96 (id=97)
$sh = $g.GroupItems.Item(94)
$sh.Left = 592.079345703125
$sh.Top = 297.3705749511719

# tx116 (id=117)
$sh = $g.GroupItems.Item(114)
$sh.Top = 536.1945190429688

# tx117 (id=118)
$sh = $g.GroupItems.Item(115)
$sh.Top = 402.61370849609375

# tx118 (id=119)
$sh = $g.GroupItems.Item(116)
$sh.Top = 269.0328369140625

# tx119 (id=120)
$sh = $g.GroupItems.Item(117)
$sh.Top = 135.4519805908203

# pl120 (id=121)
$sh = $g.GroupItems.Item(118)
$sh.Top = 539.4773559570312

# pl121 (id=122)
$sh = $g.GroupItems.Item(119)
$sh.Top = 405.896484375

# pl122 (id=123)
$sh = $g.GroupItems.Item(120)
$sh.Top = 272.31561279296875

# pl123 (id=124)
$sh = $g.GroupItems.Item(121)
$sh.Top = 138.7348175048828
